# Merge with branch evol-3577-postgresql
# The STEPS sheet gains a new "TC_STEP_CALL_DATASET" column, inserted
# right before the existing "TC_STEP_ACTION" column (G).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STEPS")

# Insert a new column before column G (TC_STEP_ACTION), shifting the
# remaining columns (TC_STEP_ACTION .. TC_STEP_CUF_<CODE>) one to the right.
$ws.Columns("G").Insert()

# Populate the header of the freshly inserted column.
$ws.Range("G1").Value = "TC_STEP_CALL_DATASET"

# Restore the active selection to sit on the new column's first data row.
$ws.Range("G2").Select() | Out-Null
